$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3399
$ws.Range("I40").Value = 2748.75
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 2748.75
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -2573.75
$ws.Range("N40").Value = -6350
$ws.Range("H43").Value = 9779.091
$ws.Range("I43").Value = 20500
$ws.Range("J43").Value = 5758.75
$ws.Range("K43").Value = 20500
$ws.Range("L43").Value = 5758.75
$ws.Range("M43").Value = -20431
$ws.Range("N43").Value = -5896.75
$ws.Range("H62").Value = 8689.166999999999
$ws.Range("I62").Value = 10438.5
$ws.Range("J62").Value = 7814.5
$ws.Range("K62").Value = 10438.5
$ws.Range("L62").Value = 7814.5
$ws.Range("M62").Value = -9814.5
$ws.Range("N62").Value = -9062.5
$ws.Range("H65").Value = 8689.166999999999
$ws.Range("I65").Value = 10438.5
$ws.Range("J65").Value = 7814.5
$ws.Range("K65").Value = 52192.5
$ws.Range("L65").Value = 39072.5
$ws.Range("M65").Value = -49072.5
$ws.Range("N65").Value = -45312.5
$ws.Range("H69").Value = 200503.2
$ws.Range("I69").Value = 15000
$ws.Range("J69").Value = 229042.16
$ws.Range("K69").Value = 45000
$ws.Range("L69").Value = 687126.48
$ws.Range("M69").Value = -44126
$ws.Range("N69").Value = -688874.48
$ws.Range("H70").Value = 11368.75
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 11368.75
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").Value = 34106.25
$ws.Range("N70").Value = -34646.25
$ws.Range("H72").Value = 200503.2
$ws.Range("I72").Value = 15000
$ws.Range("J72").Value = 229042.16
$ws.Range("K72").Value = 135000
$ws.Range("L72").Value = 2061379.44
$ws.Range("M72").Value = -130632
$ws.Range("N72").Value = -2070115.44
$ws.Range("H73").Value = 11368.75
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 11368.75
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").Value = 34106.25
$ws.Range("N73").Value = -35978.25
$ws.Range("H87").Value = 127474.5
$ws.Range("J87").Value = 159949
$ws.Range("L87").Value = 159949
$ws.Range("N87").Value = -162445
$ws.Range("H90").Value = 127474.5
$ws.Range("J90").Value = 159949
$ws.Range("L90").Value = 479847
$ws.Range("N90").Value = -492327
$ws.Range("H106").Value = 2018
$ws.Range("I106").Value = 1874
$ws.Range("K106").Value = 1874
$ws.Range("M106").Value = -1243
$ws.Range("H111").Value = 953.4666999999999
$ws.Range("J111").Value = 1297.7142
$ws.Range("L111").Value = 3893.1426
$ws.Range("N111").Value = -10027.1426
$ws.Range("H113").Value = 6190.2856
$ws.Range("H125").Value = 1671.1428
$ws.Range("I125").Value = 1345.091
$ws.Range("K125").Value = 12105.819
$ws.Range("M125").Value = -9645.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 6306.3335
$ws.Range("I38").Value = 10019
$ws.Range("J38").Value = 4450
$ws.Range("K38").Value = 10019
$ws.Range("L38").Value = 4450
$ws.Range("M38").Value = -9552
$ws.Range("N38").Value = -5384
$ws.Range("H43").Value = 92377
$ws.Range("J43").Value = 92377
$ws.Range("L43").Value = 92377
$ws.Range("N43").Value = -93003
$ws.Range("H61").Value = 1866.8182
$ws.Range("I61").Value = 1354.6786
$ws.Range("K61").Value = 1354.6786
$ws.Range("M61").Value = -1142.6786
$ws.Range("H122").Value = 1926.3462
$ws.Range("I122").Value = 1763.1818
$ws.Range("K122").Value = 5289.5454
$ws.Range("M122").Value = -2839.5454
$ws.Range("H133").Value = 89997
$ws.Range("J133").Value = 89997
$ws.Range("L133").Value = 89997
$ws.Range("N133").Value = -95057
$ws.Range("H136").Value = 1866.8182
$ws.Range("I136").Value = 1354.6786
$ws.Range("K136").Value = 4064.0358
$ws.Range("M136").Value = -1514.0358

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 602.9545000000001
$ws.Range("I94").Value = 603
$ws.Range("K94").Value = 603
$ws.Range("M94").Value = -152
$ws.Range("H128").Value = 3360
$ws.Range("I128").Value = 3360
$ws.Range("K128").Value = 10080
$ws.Range("M128").Value = -7590
$ws.Range("H134").Value = 2367.7585
$ws.Range("I134").Value = 2166.6
$ws.Range("J134").Value = 3625
$ws.Range("K134").Value = 6499.799999999999
$ws.Range("L134").Value = 10875
$ws.Range("M134").Value = -3964.799999999999
$ws.Range("N134").Value = -15945

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11708.903
$ws.Range("I58").Value = 1177.8572
$ws.Range("J58").Value = 109998.664
$ws.Range("K58").Value = 1177.8572
$ws.Range("L58").Value = 109998.664
$ws.Range("M58").Value = -974.8571999999999
$ws.Range("N58").Value = -110404.664
$ws.Range("H76").Value = 8297.5
$ws.Range("I76").Value = 8297.5
$ws.Range("K76").Value = 8297.5
$ws.Range("M76").Value = -7982.5
$ws.Range("H79").Value = 8297.5
$ws.Range("I79").Value = 8297.5
$ws.Range("K79").Value = 8297.5
$ws.Range("M79").Value = -7205.5
$ws.Range("H136").Value = 11708.903
$ws.Range("I136").Value = 1177.8572
$ws.Range("J136").Value = 109998.664
$ws.Range("K136").Value = 3533.5716
$ws.Range("L136").Value = 329995.992
$ws.Range("M136").Value = -983.5715999999998
$ws.Range("N136").Value = -335095.992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1087.7142
$ws.Range("I5").Value = 1087.7142
$ws.Range("K5").Value = 3263.1426
$ws.Range("M5").Value = -3151.1426
$ws.Range("H29").Value = 336.33334
$ws.Range("I29").Value = 265.7143
$ws.Range("K29").Value = 797.1428999999999
$ws.Range("M29").Value = -520.1428999999999
$ws.Range("H33").Value = 89.666664
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H40").Value = 145.875
$ws.Range("I40").Value = 14.8
$ws.Range("J40").Value = 364.33334
$ws.Range("K40").Value = 59.2
$ws.Range("L40").Value = 1457.33336
$ws.Range("M40").Value = 9.799999999999997
$ws.Range("N40").Value = -1595.33336
$ws.Range("H44").Value = 3000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H131").Value = 1929.3636
$ws.Range("I131").Value = 1588.5714
$ws.Range("J131").Value = 2525.75
$ws.Range("K131").Value = 4765.7142
$ws.Range("L131").Value = 7577.25
$ws.Range("M131").Value = 274.2857999999997
$ws.Range("N131").Value = -17657.25
$ws.Range("H135").Value = 1087.7142
$ws.Range("I135").Value = 1087.7142
$ws.Range("K135").Value = 9789.427799999999
$ws.Range("M135").Value = -7254.427799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 49761
$ws.Range("J96").Value = 49761
$ws.Range("L96").Value = 49761
$ws.Range("N96").Value = -55253
$ws.Range("H102").Value = 1198.7572
$ws.Range("I102").Value = 1114.0317
$ws.Range("K102").Value = 1114.0317
$ws.Range("M102").Value = 507.9683
$ws.Range("H113").Value = 1011
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2984
$ws.Range("I122").Value = 2454.2666
$ws.Range("J122").Value = 4308.3335
$ws.Range("K122").Value = 7362.7998
$ws.Range("L122").Value = 12925.0005
$ws.Range("M122").Value = -4912.7998
$ws.Range("N122").Value = -17825.0005
$ws.Range("H132").Value = 4705.7646
$ws.Range("I132").Value = 4399.933
$ws.Range("J132").Value = 6999.5
$ws.Range("K132").Value = 13199.799
$ws.Range("L132").Value = 20998.5
$ws.Range("M132").Value = -10669.799
$ws.Range("N132").Value = -26058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1974.3636
$ws.Range("I22").Value = 1859.4286
$ws.Range("K22").Value = 1859.4286
$ws.Range("M22").Value = -1564.4286
$ws.Range("H27").Value = 1974.3636
$ws.Range("I27").Value = 1859.4286
$ws.Range("K27").Value = 1859.4286
$ws.Range("M27").Value = -1752.4286
$ws.Range("H46").Value = 5739.5
$ws.Range("I46").Value = 1167.6666
$ws.Range("J46").Value = 8482.6
$ws.Range("K46").Value = 1167.6666
$ws.Range("L46").Value = 8482.6
$ws.Range("M46").Value = -979.6666
$ws.Range("N46").Value = -8858.6
$ws.Range("H55").Value = 1432.5
$ws.Range("I55").Value = 400
$ws.Range("J55").Value = 1988.4615
$ws.Range("K55").Value = 400
$ws.Range("L55").Value = 1988.4615
$ws.Range("M55").Value = -227
$ws.Range("N55").Value = -2334.4615
$ws.Range("H82").Value = 5082.409
$ws.Range("I82").Value = 3327.6
$ws.Range("J82").Value = 8842.714
$ws.Range("K82").Value = 3327.6
$ws.Range("L82").Value = 8842.714
$ws.Range("M82").Value = -2966.6
$ws.Range("N82").Value = -9564.714
$ws.Range("H85").Value = 5082.409
$ws.Range("I85").Value = 3327.6
$ws.Range("J85").Value = 8842.714
$ws.Range("K85").Value = 3327.6
$ws.Range("L85").Value = 8842.714
$ws.Range("M85").Value = -2079.6
$ws.Range("N85").Value = -11338.714
$ws.Range("H93").Value = 3217.3333
$ws.Range("I93").Value = 2439.8
$ws.Range("K93").Value = 2439.8
$ws.Range("M93").Value = -1191.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 100000
$ws.Range("J93").Value = 100000
$ws.Range("L93").Value = 100000
$ws.Range("N93").Value = -104992
$ws.Range("H113").Value = 1905.2222
$ws.Range("I113").Value = 3049
$ws.Range("J113").Value = 1578.4286
$ws.Range("K113").Value = 9147
$ws.Range("L113").Value = 4735.2858
$ws.Range("M113").Value = -6977
$ws.Range("N113").Value = -9075.2858
$ws.Range("H122").Value = 67418.97
$ws.Range("I122").Value = 83367.36
$ws.Range("K122").Value = 250102.08
$ws.Range("M122").Value = -247652.08
